$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new change-log entry row (row 21), matching the style/format of the
# existing rows above it (date in col A, year in col B, sheet name in col C,
# change description in col D).
$ws.Range("A21").Value = 43494
$ws.Range("A21").NumberFormat = "d-mmm-yy"

$ws.Range("B21").Value = 2012

$ws.Range("C21").Value = "West Beach MID"
$ws.Range("D21").Value = 'Changed blanks for samplers and recorders to "???"'

# Update the active selection to the newly added cell, matching the edit.
$ws.Range("E21").Select()
